# Anonymize the author names / subject / titles in the "Data Example - PPKSS"
# worksheet (sheet 3) and refresh the related layout (column widths, row
# heights, selection) to match the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Example - PPKSS")
$ws.Activate()

# --- Anonymized data (row 2: Bashir citation / row 3: Nerys citation) ---
$ws.Range("D2").Value = "DEEP SPACE NINE"
$ws.Range("E2").Value = "Sisko, Benjamin (UNITED FEDERATION OF PLANETS); Sisko, Benjamin (UNITED FEDERATION OF PLANETS); Bashir, Julian (UNITED FEDERATION OF PLANETS); Bashir, Julian (UNITED FEDERATION OF PLANETS)"
$ws.Range("F2").Value = "Frontier Medicine: A report documenting the trials and tribulations of medical practice in Deep Space Nine"

$ws.Range("D3").Value = "DEEP SPACE NINE"
$ws.Range("E3").Value = "Sisko, Benjamin (UNITED FEDERATION OF PLANETS); Sisko, Benjamin (UNITED FEDERATION OF PLANETS); Nerys, Kira (UNIV OF BAJOR);"
$ws.Range("F3").Value = "Key improvements resulting from the continued presence of the Federation in Bajor"

# --- Column widths for the (now longer) Authors / Title columns ---
$ws.Columns.Item(5).ColumnWidth = 27.375
$ws.Columns.Item(6).ColumnWidth = 25

# --- Row heights recomputed for the new wrapped text ---
$ws.Rows.Item(1).RowHeight = 51
$ws.Rows.Item(2).RowHeight = 99
$ws.Rows.Item(3).RowHeight = 51

# --- Selection / view housekeeping to match the re-saved workbook ---
$excel.ActiveWindow.View = -4143
$excel.ActiveWindow.Zoom = 100
$ws.Range("F8").Select()
